# Duplicate the first sheet ("New Sheet", which holds B1="New Data") twice,
# inserting both copies immediately after it (before "Introduction_Modified"),
# then rename everything to match the target layout:
#
#   1. New Sheet1             (original "New Sheet", content unchanged)
#   2. New Sheet               (copy of "New Sheet")
#   3. Introduction_Modified1  (copy of "New Sheet")
#   4. Introduction_Modified   (untouched)
#   5. Sheet                   (untouched)

$wb = $excel.ActiveWorkbook

$orig = $wb.Worksheets.Item("New Sheet")

# First duplicate, placed right after the original.
$orig.Copy($null, $orig)
$copy1 = $wb.Worksheets.Item(2)

# Second duplicate, placed right after the first duplicate.
$copy1.Copy($null, $copy1)
$copy2 = $wb.Worksheets.Item(3)

# Rename in an order that avoids transient name collisions.
$copy2.Name = "Introduction_Modified1"
$orig.Name = "New Sheet1"
$copy1.Name = "New Sheet"
